$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.382.58'
$ws.Range('E2').Value = '  -1.29%  '
$ws.Range('D3').Value = '3.451.27'
$ws.Range('E3').Value = '  -0.38%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '574.71'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.64'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.24%  '
$ws.Range('D7').Value = '3.455.21'
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.479'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.58%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '7.59'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.52%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.126'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.03%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.387'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.99%  '
$ws.Range('D13').Value = '4.051.04'
$ws.Range('E13').Value = '  -0.10%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.53'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.75%  '
$ws.Range('E15').Value = '  -0.39%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000174'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.56%  '
$ws.Range('D17').Value = '3.462.31'
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('D18').Value = '61.619.43'
$ws.Range('E18').Value = '  -0.77%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.38'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.64%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.51'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +2.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.44'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.59%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '397.66'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.64%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.568'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.63%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.69'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.39%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.997'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.59%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000124'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.53%  '
$ws.Range('D27').Value = '3.587.86'
$ws.Range('E27').Value = '  +0.00%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.179'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.78%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.62'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -2.80%  '
$ws.Range('E30').Value = '  +0.23%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.21'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.86%  '
$ws.Range('B32').Value = 'PancakeSwap'
$ws.Range('C32').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.18'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.03%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.45'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -8.32%  '
$ws.Range('E34').Value = '  +0.02%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.96'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.47%  '
$ws.Range('B36').Value = 'RenzoRestakedETH'
$ws.Range('C36').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D36').Value = '3.494.28'
$ws.Range('E36').Value = '  +0.28%  '
$ws.Range('B37').Value = 'Aptos'
$ws.Range('C37').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '7.04'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.52%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '5.18'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.97%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.56'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.65%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '167.44'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '28.36'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +7.49%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0791'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.29%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.805'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.03%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '4.53'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.10%  '
$ws.Range('B45').Value = 'FirstDigitalUSD'
$ws.Range('C45').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.13%  '
$ws.Range('B46').Value = 'Stacks'
$ws.Range('C46').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.72'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.44%  '
$ws.Range('D47').Value = '2.628.86'
$ws.Range('E47').Value = '  -0.80%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.15'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.65%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '6.96'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.67%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '22.99'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.29%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.42'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.96%  '
